$p = $ppt.ActivePresentation

$oldDate = "27/05/2015"
$newDate = "06/06/2015"

function Update-DateShape($shp) {
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# --- Slide Master: "Date Placeholder" field (27/05/2015 -> 06/06/2015) ---
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# --- Slide Layouts (all 11): same date field ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# --- Notes Master: same date field ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateShape $notesMaster.Shapes.Item($i)
}

# --- Slide 14: merge "Virus de " + "Enlace o " runs into a single run ---
$slide14 = $p.Slides.Item(14)
for ($i = 1; $i -le $slide14.Shapes.Count; $i++) {
    $shp = $slide14.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "Virus de Enlace o Directorio*") {
            $merged = $tr.Characters(1, 18)
            $merged.Text = "Virus de Enlace o "
        }
    }
}
